$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 15 (hunk 0)
$ws.Range("H15").Value = 1097.8158
$ws.Range("I15").Value = 1097.8158
$ws.Range("K15").Value = 3293.4474
$ws.Range("M15").Value = -3124.4474
# row 17 (hunk 1)
$ws.Range("H17").Value = 2443
$ws.Range("J17").Value = 2443
$ws.Range("L17").Value = 7329
$ws.Range("N17").Value = -7665
# row 88 (hunk 2)
$ws.Range("H88").Value = 1946.4546
$ws.Range("I88").Value = 1150.5
$ws.Range("K88").Value = 1150.5
$ws.Range("M88").Value = -744.5
# row 91 (hunk 3)
$ws.Range("H91").Value = 1946.4546
$ws.Range("I91").Value = 1150.5
$ws.Range("K91").Value = 1150.5
$ws.Range("M91").Value = 253.5
# row 137 (hunk 4)
$ws.Range("H137").Value = 4567.7856
$ws.Range("I137").Value = 4555.92
$ws.Range("J137").Value = 4666.6665
$ws.Range("K137").Value = 13667.76
$ws.Range("L137").Value = 13999.9995
$ws.Range("M137").Value = -11117.76
$ws.Range("N137").Value = -19099.9995
# row 138 (hunk 5)
$ws.Range("H138").Value = 6167.875
$ws.Range("I138").Value = 2355.889
$ws.Range("J138").Value = 7659.522
$ws.Range("K138").Value = 7067.667
$ws.Range("L138").Value = 22978.566
$ws.Range("M138").Value = -1927.667
$ws.Range("N138").Value = -33258.566

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32 (hunk 6)
$ws.Range("H32").Value = 3294.199
$ws.Range("I32").Value = 3013.711
$ws.Range("J32").Value = 9359.75
$ws.Range("K32").Value = 3013.711
$ws.Range("L32").Value = 9359.75
$ws.Range("M32").Value = -2726.711
$ws.Range("N32").Value = -9933.75
# row 44 (hunk 7)
$ws.Range("H44").Value = 45000
$ws.Range("J44").Value = 45000
$ws.Range("L44").Value = 45000
$ws.Range("N44").Value = -45976
# row 55 (hunk 8)
$ws.Range("H55").Value = 38666.668
$ws.Range("I55").Value = 20000
$ws.Range("J55").Value = 48000
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 48000
$ws.Range("M55").Value = -19685
$ws.Range("N55").Value = -48630
# row 74 (hunk 9)
$ws.Range("H74").Value = 1134.878
$ws.Range("I74").Value = 1059.9656
$ws.Range("J74").Value = 1315.9166
$ws.Range("K74").Value = 1059.9656
$ws.Range("L74").Value = 1315.9166
$ws.Range("M74").Value = -185.9656
$ws.Range("N74").Value = -3063.9166
# row 77 (hunk 10)
$ws.Range("H77").Value = 1134.878
$ws.Range("I77").Value = 1059.9656
$ws.Range("J77").Value = 1315.9166
$ws.Range("K77").Value = 5299.828
$ws.Range("L77").Value = 6579.583000000001
$ws.Range("M77").Value = -931.8279999999995
$ws.Range("N77").Value = -15315.583
# row 132 (hunk 11)
$ws.Range("H132").Value = 2624.762
$ws.Range("I132").Value = 2753.7632
$ws.Range("K132").Value = 8261.2896
$ws.Range("M132").Value = -5731.2896

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 86 (hunk 12)
$ws.Range("H86").Value = 2127169.5
$ws.Range("I86").Value = 2835192.5
$ws.Range("J86").Value = 3100
$ws.Range("K86").Value = 2835192.5
$ws.Range("L86").Value = 3100
$ws.Range("M86").Value = -2834069.5
$ws.Range("N86").Value = -5346
# row 89 (hunk 13)
$ws.Range("H89").Value = 2127169.5
$ws.Range("I89").Value = 2835192.5
$ws.Range("J89").Value = 3100
$ws.Range("K89").Value = 14175962.5
$ws.Range("L89").Value = 15500
$ws.Range("M89").Value = -14170346.5
$ws.Range("N89").Value = -26732

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 16 (hunk 14)
$ws.Range("H16").Value = 2486.7778
$ws.Range("I16").Value = 2740.2856
$ws.Range("J16").Value = 1599.5
$ws.Range("K16").Value = 2740.2856
$ws.Range("L16").Value = 1599.5
$ws.Range("M16").Value = -2453.2856
$ws.Range("N16").Value = -2173.5
# row 31 (hunk 15)
$ws.Range("H31").Value = 36321.613
$ws.Range("I31").Value = 2459
$ws.Range("J31").Value = 42833.652
$ws.Range("K31").Value = 2459
$ws.Range("L31").Value = 42833.652
$ws.Range("M31").Value = -2164
$ws.Range("N31").Value = -43423.652
# row 34 (hunk 16)
$ws.Range("H34").Value = 36321.613
$ws.Range("I34").Value = 2459
$ws.Range("J34").Value = 42833.652
$ws.Range("K34").Value = 2459
$ws.Range("L34").Value = 42833.652
$ws.Range("M34").Value = -2257
$ws.Range("N34").Value = -43237.652
# row 50 (hunk 17)
$ws.Range("H50").Value = 36589.125
$ws.Range("I50").Value = 24082.6
$ws.Range("J50").Value = 57433.332
$ws.Range("K50").Value = 24082.6
$ws.Range("L50").Value = 57433.332
$ws.Range("M50").Value = -23457.6
$ws.Range("N50").Value = -58683.332
# row 59 (hunk 18)
$ws.Range("H59").Value = 31175.572
$ws.Range("I59").Value = 50000
$ws.Range("J59").Value = 28038.166
$ws.Range("K59").Value = 50000
$ws.Range("L59").Value = 28038.166
$ws.Range("M59").Value = -48855
$ws.Range("N59").Value = -30328.166
# row 60 (hunk 19)
$ws.Range("H60").Value = 42031
$ws.Range("I60").Value = 6093
$ws.Range("J60").Value = 60000
$ws.Range("K60").Value = 6093
$ws.Range("L60").Value = 60000
$ws.Range("M60").Value = -5582
$ws.Range("N60").Value = -61022
# row 99 (hunk 20)
$ws.Range("H99").Value = 4637.7
$ws.Range("I99").Value = 4653
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 4653
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -3155
$ws.Range("N99").Value = -7496
# row 113 (hunk 21)
$ws.Range("H113").Value = 2486.7778
$ws.Range("I113").Value = 2740.2856
$ws.Range("J113").Value = 1599.5
$ws.Range("K113").Value = 2740.2856
$ws.Range("L113").Value = 1599.5
$ws.Range("M113").Value = -570.2856000000002
$ws.Range("N113").Value = -5939.5
# row 126 (hunk 22)
$ws.Range("H126").Value = 4637.7
$ws.Range("I126").Value = 4653
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 13959
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -11489
$ws.Range("N126").Value = -18440
# row 132 (hunk 23)
$ws.Range("H132").Value = 1553.091
$ws.Range("I132").Value = 1481.3684
$ws.Range("K132").Value = 4444.1052
$ws.Range("M132").Value = -1914.1052

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 107 (hunk 24)
$ws.Range("H107").Value = 59339.39
$ws.Range("I107").Value = 1060.9
$ws.Range("K107").Value = 3182.7
$ws.Range("M107").Value = -1262.7

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 122 (hunk 25)
$ws.Range("H122").Value = 4277.75
$ws.Range("I122").Value = 2762.125
$ws.Range("K122").Value = 8286.375
$ws.Range("M122").Value = -5836.375
# row 126 (hunk 26)
$ws.Range("H126").Value = 3647.92
$ws.Range("I126").Value = 3336.6428
$ws.Range("J126").Value = 4044.0908
$ws.Range("K126").Value = 10009.9284
$ws.Range("L126").Value = 12132.2724
$ws.Range("M126").Value = -7539.928400000001
$ws.Range("N126").Value = -17072.2724
# row 128 (hunk 27)
$ws.Range("H128").Value = 73998
$ws.Range("J128").Value = 73998
$ws.Range("L128").Value = 73998
$ws.Range("N128").Value = -83958
# row 132 (hunk 28)
$ws.Range("H132").Value = 101538.45
$ws.Range("I132").Value = 10241
$ws.Range("K132").Value = 30723
$ws.Range("M132").Value = -28193

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7 (hunk 29)
$ws.Range("H7").Value = 6947.7144
$ws.Range("I7").Value = 8026.3
$ws.Range("J7").Value = 4251.25
$ws.Range("K7").Value = 8026.3
$ws.Range("L7").Value = 4251.25
$ws.Range("M7").Value = -7914.3
$ws.Range("N7").Value = -4475.25
# row 9 (hunk 30)
$ws.Range("H9").Value = 1882.75
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 3755.5
$ws.Range("K9").Value = 10
$ws.Range("L9").Value = 3755.5
$ws.Range("M9").Value = 214
$ws.Range("N9").Value = -4203.5
# row 22 (hunk 31)
$ws.Range("H22").Value = 591.125
$ws.Range("I22").Value = 633.75
$ws.Range("J22").Value = 548.5
$ws.Range("K22").Value = 633.75
$ws.Range("L22").Value = 548.5
$ws.Range("M22").Value = -338.75
$ws.Range("N22").Value = -1138.5
# row 27 (hunk 32)
$ws.Range("H27").Value = 591.125
$ws.Range("I27").Value = 633.75
$ws.Range("J27").Value = 548.5
$ws.Range("K27").Value = 633.75
$ws.Range("L27").Value = 548.5
$ws.Range("M27").Value = -526.75
$ws.Range("N27").Value = -762.5
# row 40 (hunk 33)
$ws.Range("H40").Value = 3609.4546
$ws.Range("I40").Value = 3220.4
$ws.Range("K40").Value = 3220.4
$ws.Range("M40").Value = -3084.4
# row 55 (hunk 34)
$ws.Range("H55").Value = 1116.0476
$ws.Range("I55").Value = 406.73334
$ws.Range("K55").Value = 406.73334
$ws.Range("M55").Value = -233.73334
# row 126 (hunk 35)
$ws.Range("H126").Value = 6947.7144
$ws.Range("I126").Value = 8026.3
$ws.Range("J126").Value = 4251.25
$ws.Range("K126").Value = 24078.9
$ws.Range("L126").Value = 12753.75
$ws.Range("M126").Value = -21608.9
$ws.Range("N126").Value = -17693.75
# row 132 (hunk 36)
$ws.Range("H132").Value = 6552.85
$ws.Range("I132").Value = 5402.1333
$ws.Range("J132").Value = 10005
$ws.Range("K132").Value = 16206.3999
$ws.Range("L132").Value = 30015
$ws.Range("M132").Value = -13676.3999
$ws.Range("N132").Value = -35075
# row 136 (hunk 37)
$ws.Range("H136").Value = 403199.56
$ws.Range("I136").Value = 628648.6
$ws.Range("K136").Value = 1885945.8
$ws.Range("M136").Value = -1883395.8

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 126 (hunk 38)
$ws.Range("H126").Value = 1226.75
$ws.Range("I126").Value = 1302.3334
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 3907.0002
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -1437.0002
$ws.Range("N126").Value = -7940
# row 132 (hunk 39)
$ws.Range("H132").Value = 17602.861
$ws.Range("I132").Value = 1590.08
$ws.Range("K132").Value = 4770.24
$ws.Range("M132").Value = -2240.24
# row 136 (hunk 40)
$ws.Range("H136").Value = 16466746
$ws.Range("I136").Value = 21486792
$ws.Range("J136").Value = 402599.4
$ws.Range("K136").Value = 64460376
$ws.Range("L136").Value = 1207798.2
$ws.Range("M136").Value = -64457826
$ws.Range("N136").Value = -1212898.2
